$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B4 cell value (was 1, now 0.9)
$ws.Range("B4").Value = 0.9

# Update the active selection to B4 (was C5)
$ws.Range("B4").Select()
